$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed "Actual Consumption (MW)" / "Timestamp" series: the existing
# 28 data rows move to the new fetch window and 13 new rows (29-41) are
# appended for the added 3D Steel load in the portfolio forecast.
# Flat triples: row, consumption (col A), timestamp serial (col B).
$data = @(
    2, 5340, 45889,
    3, 5303, 45889.01041666666,
    4, 5242, 45889.02083333334,
    5, 5233, 45889.03125,
    6, 5237, 45889.04166666666,
    7, 5215, 45889.05208333334,
    8, 5164, 45889.0625,
    9, 5178, 45889.07291666666,
    10, 5196, 45889.08333333334,
    11, 5172, 45889.09375,
    12, 5140, 45889.10416666666,
    13, 5133, 45889.11458333334,
    14, 5093, 45889.125,
    15, 5055, 45889.13541666666,
    16, 5070, 45889.14583333334,
    17, 5080, 45889.15625,
    18, 5142, 45889.16666666666,
    19, 5190, 45889.17708333334,
    20, 5229, 45889.1875,
    21, 5260, 45889.19791666666,
    22, 5356, 45889.20833333334,
    23, 5404, 45889.21875,
    24, 5507, 45889.22916666666,
    25, 5578, 45889.23958333334,
    26, 5725, 45889.25,
    27, 5828, 45889.26041666666,
    28, 5881, 45889.27083333334,
    29, 5920, 45889.28125,
    30, 5948, 45889.29166666666,
    31, 5915, 45889.30208333334,
    32, 5890, 45889.3125,
    33, 5738, 45889.32291666666,
    34, 5685, 45889.33333333334,
    35, 5647, 45889.34375,
    36, 5584, 45889.35416666666,
    37, 5463, 45889.36458333334,
    38, 5347, 45889.375,
    39, 5230, 45889.38541666666,
    40, 5181, 45889.39583333334,
    41, 5148, 45889.40625
)

for ($i = 0; $i -lt $data.Length; $i += 3) {
    $r = $data[$i]
    $a = $data[$i + 1]
    $b = $data[$i + 2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
